# [Feat 2269] moar debug
# Update the "call step" path labels in column G to be prefixed with
# "CALL " (e.g. "call/path/1" -> "CALL /path/1") and restore the
# previously-selected active cell to G4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STEPS")

$ws.Range("G2").Value = "CALL /path/1"
$ws.Range("G3").Value = "CALL /path/2"
$ws.Range("G4").Value = "CALL /path/3"

$ws.Activate() | Out-Null
$ws.Range("G4").Select() | Out-Null
